$d = $word.ActiveDocument

# Locate the paragraph whose entire text is "Update 16" (the last of the
# "Update N" heading paragraphs at the top of the document) so we can add
# a new "Update 17" paragraph right after it, matching its formatting.
$target = $null
foreach ($p in $d.Paragraphs) {
    $text = $p.Range.Text
    $text = $text.TrimEnd([char]13, [char]7)
    if ($text -eq "Update 16") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $r = $target.Range
    $r.Collapse(0)            # wdCollapseEnd -> move to just before the paragraph mark
    $r.InsertAfter("Update 17" + [char]13)
}
